# Update cryptocurrency price/volume figures per the Feb 29 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.131.85"
$ws.Range("E2").Value = "  +3.17%  "
$ws.Range("D3").Value = "3.408.73"
$ws.Range("E3").Value = "  +3.81%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "406.94"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.00"
$ws.Range("E6").Value = "  +17.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.608"
$ws.Range("E7").Value = "  +7.44%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +9.52%  "
$ws.Range("E10").Value = "  +12.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.17"
$ws.Range("E11").Value = "  +9.19%  "
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("D13").Value = "3.963.66"
$ws.Range("E13").Value = "  +4.45%  "
$ws.Range("E14").Value = "  +5.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.79"
$ws.Range("E15").Value = "  +4.51%  "
$ws.Range("D16").Value = "3.420.61"
$ws.Range("E16").Value = "  +2.49%  "
$ws.Range("D17").Value = "62.056.50"
$ws.Range("E17").Value = "  +3.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.52"
$ws.Range("E18").Value = "  +8.75%  "
$ws.Range("E19").Value = "  +5.32%  "
$ws.Range("E20").Value = "  +17.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.28"
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "82.79"
$ws.Range("E22").Value = "  +12.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.18"
$ws.Range("E23").Value = "  +6.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "308.65"
$ws.Range("E24").Value = "  +4.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.17"
$ws.Range("E25").Value = "  +3.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.59"
$ws.Range("E26").Value = "  +15.05%  "
$ws.Range("E27").Value = "  +11.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "29.79"
$ws.Range("E28").Value = "  +2.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.174"
$ws.Range("E29").Value = "  +1.85%  "
$ws.Range("E30").Value = "  +1.30%  "
$ws.Range("E31").Value = "  +2.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.79"
$ws.Range("E32").Value = "  +5.70%  "
$ws.Range("E33").Value = "  +6.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "42.54"
$ws.Range("E34").Value = "  +9.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E36").Value = "  +2.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.43"
$ws.Range("E37").Value = "  +0.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.43"
$ws.Range("E39").Value = "  +4.09%  "
$ws.Range("E40").Value = "  -3.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.03"
$ws.Range("E41").Value = "  +9.17%  "
$ws.Range("E42").Value = "  +5.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "137.60"
$ws.Range("E43").Value = "  +2.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.98"
$ws.Range("E44").Value = "  +5.66%  "
$ws.Range("E45").Value = "  -2.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.07"
$ws.Range("E46").Value = "  +5.36%  "
$ws.Range("E47").Value = "  +2.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.71"
$ws.Range("E48").Value = "  +4.30%  "
$ws.Range("D49").Value = "3.750.20"
$ws.Range("E49").Value = "  +4.17%  "
$ws.Range("D50").Value = "2.149.81"
$ws.Range("E50").Value = "  +1.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.35"
$ws.Range("E51").Value = "  -0.69%  "
